# Auto-generated script applying the Moogle_Profits market-price/profit refresh
# Updates H/I/J/K/L/M/N columns (price & profit data) across the ALC, ARM, BSM, CRP,
# CUL, GSM, LTW and WVR sheets; clears a couple of now-empty profit cells; adds a few
# previously-missing profit cells.
$wb = $excel.ActiveWorkbook

### Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 50000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 50000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H51").Value = 8006.25
$ws.Range("J51").Value = 8207.5
$ws.Range("L51").Value = 8207.5
$ws.Range("N51").Value = -9175.5
$ws.Range("H69").Value = 18094.38
$ws.Range("J69").Value = 18631.158
$ws.Range("L69").Value = 55893.474
$ws.Range("N69").Value = -57641.474
$ws.Range("H72").Value = 18094.38
$ws.Range("J72").Value = 18631.158
$ws.Range("L72").Value = 167680.422
$ws.Range("N72").Value = -176416.422
$ws.Range("H86").Value = 5423.647
$ws.Range("I86").Value = 1432.8334
$ws.Range("J86").Value = 7600.4546
$ws.Range("K86").Value = 1432.8334
$ws.Range("L86").Value = 7600.4546
$ws.Range("M86").Value = -309.8334
$ws.Range("N86").Value = -9846.454600000001
$ws.Range("H89").Value = 5423.647
$ws.Range("I89").Value = 1432.8334
$ws.Range("J89").Value = 7600.4546
$ws.Range("K89").Value = 7164.166999999999
$ws.Range("L89").Value = 38002.273
$ws.Range("M89").Value = -1548.166999999999
$ws.Range("N89").Value = -49234.273
$ws.Range("H106").Value = 15716896
$ws.Range("I106").Value = 18335824
$ws.Range("K106").Value = 18335824
$ws.Range("M106").Value = -18335193
$ws.Range("H130").Value = 142225
$ws.Range("J130").Value = 142225
$ws.Range("L130").Value = 142225
$ws.Range("N130").Value = -152265
$ws.Range("H138").Value = 3737.3684
$ws.Range("J138").Value = 6776.3
$ws.Range("L138").Value = 20328.9
$ws.Range("N138").Value = -30608.9

### Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5151.375
$ws.Range("I32").Value = 2013.7313
$ws.Range("K32").Value = 2013.7313
$ws.Range("M32").Value = -1726.7313

### Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 52074.5
$ws.Range("J21").Value = 52074.5
$ws.Range("L21").Value = 52074.5
$ws.Range("N21").Value = -52546.5
$ws.Range("H36").Value = 995
$ws.Range("I36").Value = 995
$ws.Range("K36").Value = 995
$ws.Range("M36").Value = -461
$ws.Range("H54").Value = 32845
$ws.Range("J54").Value = 45496.75
$ws.Range("L54").Value = 45496.75
$ws.Range("N54").Value = -46464.75
$ws.Range("H94").Value = 606.11536
$ws.Range("I94").Value = 606.11536
$ws.Range("K94").Value = 606.11536
$ws.Range("M94").Value = -155.11536
$ws.Range("H128").Value = 20335
$ws.Range("I128").Value = 20335
$ws.Range("K128").Value = 61005
$ws.Range("M128").Value = -58515

### Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 9001
$ws.Range("I2").Value = 8004
$ws.Range("K2").Value = 8004
$ws.Range("M2").Value = -7891
$ws.Range("H31").Value = 9747.947
$ws.Range("I31").Value = 4795.125
$ws.Range("K31").Value = 4795.125
$ws.Range("M31").Value = -4500.125
$ws.Range("H34").Value = 9747.947
$ws.Range("I34").Value = 4795.125
$ws.Range("K34").Value = 4795.125
$ws.Range("M34").Value = -4593.125
$ws.Range("H35").Value = 15749.5
$ws.Range("I35").Value = 1500
$ws.Range("J35").Value = 29999
$ws.Range("K35").Value = 1500
$ws.Range("L35").Value = 29999
$ws.Range("M35").Value = -1206
$ws.Range("N35").Value = -30587
$ws.Range("H37").Value = 30500
$ws.Range("I37").Value = 4000
$ws.Range("J37").Value = 57000
$ws.Range("K37").Value = 4000
$ws.Range("L37").Value = 57000
$ws.Range("M37").Value = -3893
$ws.Range("N37").Value = -57214
$ws.Range("H38").Value = 21999.334
$ws.Range("I38").Value = 5999
$ws.Range("K38").Value = 5999
$ws.Range("M38").Value = -5622
$ws.Range("H46").Value = 21999.334
$ws.Range("I46").Value = 5999
$ws.Range("K46").Value = 5999
$ws.Range("M46").Value = -5788
$ws.Range("H50").Value = 47816
$ws.Range("J50").Value = 66723.25
$ws.Range("L50").Value = 66723.25
$ws.Range("N50").Value = -67973.25
$ws.Range("H99").Value = 2249.6365
$ws.Range("I99").Value = 1869.4117
$ws.Range("J99").Value = 3542.4
$ws.Range("K99").Value = 1869.4117
$ws.Range("L99").Value = 3542.4
$ws.Range("M99").Value = -371.4117000000001
$ws.Range("N99").Value = -6538.4
$ws.Range("H122").Value = 1937.4546
$ws.Range("I122").Value = 1951.2
$ws.Range("K122").Value = 5853.6
$ws.Range("M122").Value = -3403.6
$ws.Range("H126").Value = 2249.6365
$ws.Range("I126").Value = 1869.4117
$ws.Range("J126").Value = 3542.4
$ws.Range("K126").Value = 5608.2351
$ws.Range("L126").Value = 10627.2
$ws.Range("M126").Value = -3138.2351
$ws.Range("N126").Value = -15567.2
$ws.Range("H134").Value = 5653.636
$ws.Range("I134").Value = 2961.875
$ws.Range("K134").Value = 8885.625
$ws.Range("M134").Value = -6350.625

### Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 974.75
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 949.5
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 8545.5
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -13445.5

### Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4288.017
$ws.Range("I70").Value = 2385.111
$ws.Range("K70").Value = 2385.111
$ws.Range("M70").Value = -2115.111
$ws.Range("H73").Value = 4288.017
$ws.Range("I73").Value = 2385.111
$ws.Range("K73").Value = 2385.111
$ws.Range("M73").Value = -1449.111
$ws.Range("H80").Value = 6650.5757
$ws.Range("I80").Value = 4915.8237
$ws.Range("K80").Value = 4915.8237
$ws.Range("M80").Value = -3917.8237
$ws.Range("H83").Value = 6650.5757
$ws.Range("I83").Value = 4915.8237
$ws.Range("K83").Value = 24579.1185
$ws.Range("M83").Value = -19587.1185
$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080

### Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 47342.68
$ws.Range("I7").Value = 57091.9
$ws.Range("J7").Value = 8345.799999999999
$ws.Range("K7").Value = 57091.9
$ws.Range("L7").Value = 8345.799999999999
$ws.Range("M7").Value = -56979.9
$ws.Range("N7").Value = -8569.799999999999
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H26").Value = 56112.11
$ws.Range("I26").Value = 25004.5
$ws.Range("K26").Value = 25004.5
$ws.Range("M26").Value = -24709.5
$ws.Range("H55").Value = 1415.6428
$ws.Range("I55").Value = 771.6667
$ws.Range("K55").Value = 771.6667
$ws.Range("M55").Value = -598.6667
$ws.Range("H93").Value = 2152.6316
$ws.Range("I93").Value = 2000.6923
$ws.Range("K93").Value = 2000.6923
$ws.Range("M93").Value = -752.6922999999999
$ws.Range("H126").Value = 47342.68
$ws.Range("I126").Value = 57091.9
$ws.Range("J126").Value = 8345.799999999999
$ws.Range("K126").Value = 171275.7
$ws.Range("L126").Value = 25037.4
$ws.Range("M126").Value = -168805.7
$ws.Range("N126").Value = -29977.4

### Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 27001
$ws.Range("I17").Value = 14943.077
$ws.Range("J17").Value = 105377.5
$ws.Range("K17").Value = 14943.077
$ws.Range("L17").Value = 105377.5
$ws.Range("M17").Value = -14771.077
$ws.Range("N17").Value = -105721.5
$ws.Range("H41").Value = 17018.25
$ws.Range("J41").Value = 16991
$ws.Range("L41").Value = 16991
$ws.Range("N41").Value = -17771
$ws.Range("H122").Value = 4594
$ws.Range("I122").Value = 3464.8
$ws.Range("K122").Value = 10394.4
$ws.Range("M122").Value = -7944.400000000001

Write-Host "Applied Moogle_Profits scheduled-runner update across all sheets"
